$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the fruit name in row 3 from "Apple" to "Murgi"
$ws.Range("B3").Value = "Murgi"
